# Component_Selection.xlsx edit script
# "Schematic cleaned up, connectors labeled"
#
# Applies the Shunt Resistors sheet recalculation (new shunt values +
# thermal-fusing check formulas) and restores the two sheets' scroll/
# selection state.

$wb = $excel.ActiveWorkbook

$wsShunt = $wb.Worksheets.Item("Shunt Resistors")
$wsMosfet = $wb.Worksheets.Item("MOSFETs")

# --- Shunt Resistors: "Continuous Shunt" block (rows 3-7) -------------
# Continuous current bumped 60A -> 200A, continuous shunt value 0.01 -> 1E-3
$wsShunt.Range("B3").Value = 200
$wsShunt.Range("B4").Value = 0.001

# Add "time to fuse" helper formulas (column H) for each candidate shunt,
# formatted to match the scientific notation already used in column F.
$wsShunt.Range("H4").Formula = "=60*60*F4"
$wsShunt.Range("H4").NumberFormat = "##0.0E+0"

$wsShunt.Range("H5").Formula = "=30*30*F5"
$wsShunt.Range("H5").NumberFormat = "##0.0E+0"

$wsShunt.Range("H6").Formula = "=16*16*F6"
$wsShunt.Range("H6").NumberFormat = "##0.0E+0"

$wsShunt.Range("H7").Formula = "=5*5*F7"
$wsShunt.Range("H7").NumberFormat = "##0.0E+0"

# --- Shunt Resistors: "Pulse Shunt" block (rows 19-32) -----------------
# Pulse shunt resistance 5E-4 -> 2E-4
$wsShunt.Range("B21").Value = 0.0002

# --- Restore saved selection / scroll position for both sheets --------
# Select the (currently) non-active sheet's range first so its selection
# is recorded without stealing the active tab away from MOSFETs.
$wsShunt.Range("B22").Select()

# MOSFETs is the workbook's active tab; select it last so it remains
# active/tabSelected after this script runs.
$wsMosfet.Activate()
$wsMosfet.Range("E30").Select()
